$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append three new rows to the "Logs" sheet (rows 8-10) ---

# Row 8
$logs.Range("A8").Value = "Afmelding nieuwsbrief"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D8").Value = "Afmelding"
$logs.Range("F8").Value = "2025-06-17 11:59:01"
$logs.Range("G8").Value = "Nee"

# Row 9
$logs.Range("A9").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D9").Value = "Informatieaanvraag"
$logs.Range("E9").Value = "Beste klant,`nDank voor uw interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur. Op zaterdag zijn wij geopend van 10:00 tot 16:00 uur. Op zondag zijn wij gesloten. Mocht u verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Range("F9").Value = "2025-06-17 11:59:03"
$logs.Range("G9").Value = "Ja"

# Row 10
$logs.Range("A10").Value = "Re: Wat zijn jullie openingstijden?"
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Beste klant,`nDank voor uw interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur. Op zaterdag zijn wij geopend van 10:00 tot 16:00 uur. Op zondag zijn wij gesloten. Mocht u verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Range("D10").Value = "Informatieaanvraag"
$logs.Range("E10").Value = "Beste klant,`nHartelijk dank voor uw interesse. Onze openingstijden zijn maandag t/m vrijdag van 9:00 tot 18:00 uur en zaterdag van 10:00 tot 16:00 uur. Op zondag zijn wij gesloten. Voor verdere vragen kunt u altijd contact met ons opnemen.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Range("F10").Value = "2025-06-17 12:29:15"
$logs.Range("G10").Value = "Ja"

# --- Extend the conditional-formatting ranges to cover the new rows ---
$fcCategory = $logs.Range("D2:D7").FormatConditions
for ($i = 1; $i -le $fcCategory.Count; $i++) {
    $fcCategory.Item($i).ModifyAppliesToRange($logs.Range("D2:D10"))
}

$fcAnswered = $logs.Range("G2:G7").FormatConditions
for ($i = 1; $i -le $fcAnswered.Count; $i++) {
    $fcAnswered.Item($i).ModifyAppliesToRange($logs.Range("G2:G10"))
}

# --- Update the "Dashboard" summary sheet ---
# New order/counts: Informatieaanvraag=4, Afmelding=2, Overig=1, Bestelling=1, Klacht=1
$dash.Range("A2").Value = "Informatieaanvraag"
$dash.Range("B2").Value = 4

$dash.Range("A3").Value = "Afmelding"
$dash.Range("B3").Value = 2

$dash.Range("A4").Value = "Overig"
$dash.Range("B4").Value = 1

$dash.Range("A5").Value = "Bestelling"
$dash.Range("B5").Value = 1

$dash.Range("A6").Value = "Klacht"
$dash.Range("B6").Value = 1
